$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange

# Target paragraph (3rd paragraph, level 1) currently holds three runs:
#   "Utilize " + "Asp.Net" (err="1") + " function attributes to prevent unauthorized access to certain regions and mitigate CSRF"
# The edit merges them into a single run (keeping the first run's formatting, i.e. no err="1").
$para = $tr.Paragraphs(3, 1)

# First set to a distinct placeholder so the engine actually rewrites the
# paragraph's runs (setting identical concatenated text is treated as a
# no-op and would leave the original 3 runs untouched).
$para.Text = "TEMP_PLACEHOLDER_TEXT"

$para2 = $tr.Paragraphs(3, 1)
$para2.Text = "Utilize Asp.Net function attributes to prevent unauthorized access to certain regions and mitigate CSRF"
